$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell $ws "D2" "279.10"
Set-TextCell $ws "E2" "6.51%"
Set-TextCell $ws "D3" "27.36"
Set-TextCell $ws "E3" "2.73%"
Set-TextCell $ws "D4" "4.798"
Set-TextCell $ws "E4" "2.07%"
Set-TextCell $ws "D5" "0.06350"
Set-TextCell $ws "D6" "6.942"
Set-TextCell $ws "E6" "3.61%"
Set-TextCell $ws "D7" "3.367"
Set-TextCell $ws "E7" "6.23%"
Set-TextCell $ws "D8" "0.8819"
Set-TextCell $ws "E8" "3.60%"
Set-TextCell $ws "D9" "0.9559"
Set-TextCell $ws "E9" "5.15%"
Set-TextCell $ws "D10" "0.1481"
Set-TextCell $ws "E10" "5.60%"
Set-TextCell $ws "D11" "0.05265"
Set-TextCell $ws "E11" "3.39%"
Set-TextCell $ws "D12" "0.07274"
Set-TextCell $ws "E12" "2.49%"
Set-TextCell $ws "D13" "0.03134"
Set-TextCell $ws "E13" "0.68%"
Set-TextCell $ws "D14" "0.09064"
Set-TextCell $ws "E14" "0.22%"
Set-TextCell $ws "D15" "0.001563"
Set-TextCell $ws "E15" "1.29%"
Set-TextCell $ws "D16" "0.0006252"
Set-TextCell $ws "E16" "1.04%"
Set-TextCell $ws "D17" "0.005818"
Set-TextCell $ws "E17" "-2.16%"
Set-TextCell $ws "E18" "0.43%"
Set-TextCell $ws "D19" "2.278"
Set-TextCell $ws "E19" "6.13%"
Set-TextCell $ws "D20" "0.3126"
Set-TextCell $ws "E20" "1.76%"
Set-TextCell $ws "D21" "0.1339"
Set-TextCell $ws "E21" "4.56%"
Set-TextCell $ws "D22" "3.868"
Set-TextCell $ws "E22" "-6.33%"
Set-TextCell $ws "D23" "0.04310"
Set-TextCell $ws "E23" "1.77%"
Set-TextCell $ws "D24" "0.001181"
Set-TextCell $ws "E24" "0.05%"
Set-TextCell $ws "E25" "5.56%"
Set-TextCell $ws "E26" "-0.11%"
Set-TextCell $ws "D27" "0.0001689"
Set-TextCell $ws "E27" "-12.83%"
Set-TextCell $ws "D40" "0.04091"
Set-TextCell $ws "E40" "3.45%"
Set-TextCell $ws "D41" "0.006641"
Set-TextCell $ws "E41" "58.40%"
Set-TextCell $ws "D42" "0.1164"
Set-TextCell $ws "E42" "4.72%"
Set-TextCell $ws "D43" "0.002287"
Set-TextCell $ws "E43" "10.96%"
Set-TextCell $ws "D44" "0.01249"
Set-TextCell $ws "E44" "-10.32%"
Set-TextCell $ws "D45" "0.00005222"
Set-TextCell $ws "E45" "2.06%"
Set-TextCell $ws "E47" "821.69%"
Set-TextCell $ws "D48" "0.02249"
Set-TextCell $ws "E48" "6.00%"
Set-TextCell $ws "D50" "0.0001999"
Set-TextCell $ws "E50" "-0.11%"
